# Update "想去人数" (F column) figures for the two worksheets that share
# this data table: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 53
    $ws.Range("F5").Value = 2425
    $ws.Range("F6").Value = 223
}
